$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Date serial used for all 4 new rows (45703 == 2025-02-15). Written as a raw
# number (not a DateTime) so the engine doesn't auto-mint a brand-new,
# one-off date number-format style for the cell; the correct m/d/yyyy
# display format (the existing style used by the rest of column A) is
# applied afterwards via a format-only paste.
$newDate = 45703

# Row 98: LatentView Analytics
$ws.Cells.Item(98, 1).Value = $newDate
$ws.Cells.Item(98, 2).Value = "LatentView Analytics"
$ws.Cells.Item(98, 3).Value = "Senior Data Scientist"
$ws.Cells.Item(98, 4).Value = "Market analyytics, but seattle, easy"
$ws.Cells.Item(98, 6).Value = "https://www.linkedin.com/jobs/view/4139907275/?refId=ByteString(length%3D16%2Cbytes%3D150983f7...8f6a262d)&trackingId=JQaQxdwWX0sVJpgbRDE9zw%3D%3D"

# Row 99: Tata Consultancy Services
$ws.Cells.Item(99, 1).Value = $newDate
$ws.Cells.Item(99, 2).Value = "Tata Consultancy Services"
$ws.Cells.Item(99, 3).Value = "Data Scientist"
$ws.Cells.Item(99, 4).Value = "they want SQL… consulting"
$ws.Cells.Item(99, 6).Value = "https://www.linkedin.com/jobs/view/4149289471/?refId=ByteString(length%3D16%2Cbytes%3Da77c4f9b...24e62942)&trackingId=fZD94CNeR4eK%2F1ctw%2B1IBg%3D%3D"

# Row 100: GLX ANALYTIX (Jobname typed before Entity, matching original authoring order)
$ws.Cells.Item(100, 1).Value = $newDate
$ws.Cells.Item(100, 3).Value = "Senior Data Scientist "
$ws.Cells.Item(100, 2).Value = "GLX ANALYTIX"
$ws.Cells.Item(100, 4).Value = "Denmark, personalized medicine"
$ws.Cells.Item(100, 6).Value = "https://www.linkedin.com/jobs/view/4150250844/?refId=ByteString(length%3D16%2Cbytes%3Db8b66b1f...f79b0827)&trackingId=IJSdoryuUWCZ8UQNmm2Csw%3D%3D"

# Row 101: SureCost
$ws.Cells.Item(101, 1).Value = $newDate
$ws.Cells.Item(101, 2).Value = "SureCost"
$ws.Cells.Item(101, 3).Value = "Senior Data Scientist"
$ws.Cells.Item(101, 4).Value = "St. Petersburg, Fl,  pharmacy inventory, more of MLOPS?"
$ws.Cells.Item(101, 6).Value = "https://www.linkedin.com/jobs/view/4152258208/?refId=I3ESEstkRIyQ0GNGtFx%2FYQ%3D%3D&trackingId=jyT38KXWQ46%2FKZ%2BxwTqcfQ%3D%3D"

# Apply date format (matching existing column A style - numFmtId 14 m/d/yyyy) to the new date
# cells by copying the format from the row above, so the existing style index is reused instead
# of minting a new (duplicate) cellXf.
$ws.Range("A97").Copy()
$ws.Range("A98:A101").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update the view state: scroll/freeze pane + selection, matching final saved state
$ws.Application.ActiveWindow.ScrollRow = 92
$ws.Range("F101").Select()
